$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing H column values (rows 22-26) to text "xxx.00" ---
$ws.Cells.Item(22, 8).Value = "'265.00"
$ws.Cells.Item(23, 8).Value = "'195.00"
$ws.Cells.Item(24, 8).Value = "'160.00"
$ws.Cells.Item(25, 8).Value = "'116.00"
$ws.Cells.Item(26, 8).Value = "'123.00"

# --- Update row 28: E (concepto), F (valor capital) and H (lawyer fee text) ---
$ws.Cells.Item(28, 5).Value = "PRESTAMOS"
$ws.Cells.Item(28, 6).Value = "'12.25"
$ws.Cells.Item(28, 8).Value = "'319.00"

# --- Add new row 30 ---
$ws.Cells.Item(30, 1).Value = "'1"
$ws.Cells.Item(30, 2).Value = "Peter Patricio Tene Ojeda"
$ws.Cells.Item(30, 3).Value = "'174582556"
$ws.Cells.Item(30, 4).Value = "'174582556001"
$ws.Cells.Item(30, 5).Value = "PRESTAMOS"
$ws.Cells.Item(30, 6).Value = "'45.85"
$ws.Cells.Item(30, 7).Value = 48.85
$ws.Cells.Item(30, 8).Value = "'63.00"
$ws.Cells.Item(30, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(30, 10).Value = "Lic. Alexander Javier Miranda Granero"

# --- Add new row 31 ---
$ws.Cells.Item(31, 1).Value = "'2"
$ws.Cells.Item(31, 2).Value = "Peter Patricio Tene Ojeda"
$ws.Cells.Item(31, 3).Value = "'174582556"
$ws.Cells.Item(31, 4).Value = "'174582556001"
$ws.Cells.Item(31, 5).Value = "PRESTAMOS"
$ws.Cells.Item(31, 6).Value = "'68.98"
$ws.Cells.Item(31, 7).Value = 98.65000000000001
$ws.Cells.Item(31, 8).Value = "'128.00"
$ws.Cells.Item(31, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(31, 10).Value = "Lic. Alexander Javier Miranda Granero"

# --- Add new row 32 ---
$ws.Cells.Item(32, 1).Value = "'3"
$ws.Cells.Item(32, 2).Value = "Peter Patricio Tene Ojeda"
$ws.Cells.Item(32, 3).Value = "'174582556"
$ws.Cells.Item(32, 4).Value = "'174582556001"
$ws.Cells.Item(32, 5).Value = "PRESTAMOS"
$ws.Cells.Item(32, 6).Value = "'48.59"
$ws.Cells.Item(32, 7).Value = 200.56
$ws.Cells.Item(32, 8).Value = "'260.00"
$ws.Cells.Item(32, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(32, 10).Value = "Lic. Alexander Javier Miranda Granero"

# --- Add new row 33 ---
$ws.Cells.Item(33, 1).Value = "'8"
$ws.Cells.Item(33, 2).Value = "Peter Patricio Tene Ojeda"
$ws.Cells.Item(33, 3).Value = "'174582556"
$ws.Cells.Item(33, 4).Value = "'174582556001"
$ws.Cells.Item(33, 5).Value = "PLANTILLA DE APORTES"
$ws.Cells.Item(33, 6).Value = "'156.23"
$ws.Cells.Item(33, 7).Value = 89.56
$ws.Cells.Item(33, 8).Value = "'116.00"
$ws.Cells.Item(33, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(33, 10).Value = "Lic. Alexander Javier Miranda Granero"
